$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update TruckID column (A) and loading time columns (B, C, D) to reflect
# the re-ordered / fixed assignments and new population generation results.

$ws.Range("A2").Value = 3

$ws.Range("A3").Value = 5
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 11
$ws.Range("D3").Value = 11

$ws.Range("A4").Value = 1
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 6

$ws.Range("A5").Value = 2
$ws.Range("C5").Value = 11
$ws.Range("D5").Value = 12

$ws.Range("C6").Value = 17
$ws.Range("D6").Value = 18

$ws.Range("C7").Value = 23
$ws.Range("D7").Value = 24

$ws.Range("C8").Value = 29
$ws.Range("D8").Value = 30
